# Updated cryptos list values (price/volume) per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.157.77"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").Value = "2.576.85"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'592.29"
$ws.Range("E5").Value = "  +1.31%  "
$ws.Range("D6").Value = "'144.53"
$ws.Range("E6").Value = "  -2.18%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -1.84%  "
$ws.Range("E9").Value = "  -2.14%  "
$ws.Range("D10").Value = "'5.59"
$ws.Range("E10").Value = "  -0.83%  "
$ws.Range("D12").Value = "'0.351"
$ws.Range("E12").Value = "  -1.62%  "
$ws.Range("D13").Value = "'27.19"
$ws.Range("E13").Value = "  -0.90%  "
$ws.Range("D14").Value = "3.040.30"
$ws.Range("E14").Value = "  +0.31%  "
$ws.Range("D15").Value = "63.065.91"
$ws.Range("E15").Value = "  -0.12%  "
$ws.Range("E16").Value = "  -1.02%  "
$ws.Range("D17").Value = "2.568.99"
$ws.Range("E17").Value = "  -0.57%  "
$ws.Range("D18").Value = "'11.09"
$ws.Range("E18").Value = "  -2.42%  "
$ws.Range("D19").Value = "'341.37"
$ws.Range("E19").Value = "  -0.58%  "
$ws.Range("D20").Value = "'4.34"
$ws.Range("E20").Value = "  -1.86%  "
$ws.Range("D21").Value = "'6.65"
$ws.Range("E21").Value = "  -3.61%  "
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("E23").Value = "  +3.70%  "
$ws.Range("D24").Value = "'67.81"
$ws.Range("E24").Value = "  +1.35%  "
$ws.Range("D25").Value = "'1.60"
$ws.Range("E25").Value = "  +7.61%  "
$ws.Range("D26").Value = "'1.62"
$ws.Range("E26").Value = "  -0.54%  "
$ws.Range("E27").Value = "  -3.06%  "
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("D29").Value = "'7.92"
$ws.Range("E29").Value = "  -2.99%  "
$ws.Range("D30").Value = "'8.26"
$ws.Range("E30").Value = "  -2.64%  "
$ws.Range("E31").Value = "  -2.10%  "
$ws.Range("D32").Value = "'469.61"
$ws.Range("E32").Value = "  +0.89%  "
$ws.Range("E33").Value = "  -3.23%  "
$ws.Range("D34").Value = "'1.68"
$ws.Range("E34").Value = "  +3.12%  "
$ws.Range("D35").Value = "'176.60"
$ws.Range("E35").Value = "  +0.41%  "
$ws.Range("D37").Value = "'0.396"
$ws.Range("E37").Value = "  -2.91%  "
$ws.Range("D38").Value = "'18.86"
$ws.Range("E38").Value = "  -1.89%  "
$ws.Range("D39").Value = "'4.57"
$ws.Range("E39").Value = "  +0.17%  "
$ws.Range("E41").Value = "  -3.43%  "
$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").Value = "'40.07"
$ws.Range("E42").Value = "  +1.16%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").Value = "'158.40"
$ws.Range("E43").Value = "  +4.27%  "
$ws.Range("D44").Value = "'3.70"
$ws.Range("E44").Value = "  -3.19%  "
$ws.Range("D45").Value = "'21.39"
$ws.Range("E45").Value = "  +1.70%  "
$ws.Range("D46").Value = "'0.634"
$ws.Range("E46").Value = "  +3.45%  "
$ws.Range("E47").Value = "  -1.59%  "
$ws.Range("D48").Value = "'0.0963"
$ws.Range("E48").Value = "  -1.68%  "
$ws.Range("D49").Value = "'0.0237"
$ws.Range("E49").Value = "  -0.96%  "
$ws.Range("E50").Value = "  -1.97%  "
$ws.Range("D51").Value = "'11.40"
$ws.Range("E51").Value = "  +0.12%  "
